$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.288.64"
$ws.Range("E2").Value = "  +0.95%  "

# Row 3
$ws.Range("D3").Value = "1.926.34"
$ws.Range("E3").Value = "  +1.14%  "

# Row 4
$origStyle_4 = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("D4").Style = $origStyle_4
$ws.Range("E4").Value = "  -0.34%  "

# Row 5
$origStyle_5 = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.43"
$ws.Range("D5").Style = $origStyle_5
$ws.Range("E5").Value = "  -0.07%  "

# Row 6
$origStyle_6 = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").Style = $origStyle_6
$ws.Range("E6").Value = "  -0.30%  "

# Row 7
$origStyle_7 = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4605"
$ws.Range("D7").Style = $origStyle_7
$ws.Range("E7").Value = "  +0.15%  "

# Row 8
$origStyle_8 = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3870"
$ws.Range("D8").Style = $origStyle_8
$ws.Range("E8").Value = "  +0.02%  "

# Row 9
$origStyle_9 = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.74"
$ws.Range("D9").Style = $origStyle_9
$ws.Range("E9").Value = "  -0.09%  "

# Row 10
$origStyle_10 = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07812"
$ws.Range("D10").Style = $origStyle_10
$ws.Range("E10").Value = "  -0.10%  "

# Row 11
$origStyle_11 = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9732"
$ws.Range("D11").Style = $origStyle_11
$ws.Range("E11").Value = "  -1.30%  "

# Row 12
$origStyle_12 = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.65"
$ws.Range("D12").Style = $origStyle_12
$ws.Range("E12").Value = "  +4.03%  "

# Row 13
$ws.Range("D13").Value = "1.919.36"
$ws.Range("E13").Value = "  -0.26%  "

# Row 14
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$origStyle_14 = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.768"
$ws.Range("D14").Style = $origStyle_14
$ws.Range("E14").Value = "  +0.63%  "

# Row 15
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$origStyle_15 = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.059"
$ws.Range("D15").Style = $origStyle_15
$ws.Range("E15").Value = "  +1.02%  "

# Row 16
$origStyle_16 = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07086"
$ws.Range("D16").Style = $origStyle_16
$ws.Range("E16").Value = "  +0.66%  "

# Row 17
$origStyle_17 = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "86.68"
$ws.Range("D17").Style = $origStyle_17
$ws.Range("E17").Value = "  -1.04%  "

# Row 18
$origStyle_18 = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.003"
$ws.Range("D18").Style = $origStyle_18
$ws.Range("E18").Value = "  -0.34%  "

# Row 19
$origStyle_19 = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000009721"
$ws.Range("D19").Style = $origStyle_19
$ws.Range("E19").Value = "  -1.71%  "

# Row 20
$origStyle_20 = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.01"
$ws.Range("D20").Style = $origStyle_20
$ws.Range("E20").Value = "  +0.23%  "

# Row 21
$ws.Range("E21").Value = "  +0.08%  "

# Row 22
$ws.Range("D22").Value = "29.298.61"
$ws.Range("E22").Value = "  +0.62%  "

# Row 23
$origStyle_23 = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.484"
$ws.Range("D23").Style = $origStyle_23
$ws.Range("E23").Value = "  +3.23%  "

# Row 24
$origStyle_24 = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.09"
$ws.Range("D24").Style = $origStyle_24
$ws.Range("E24").Value = "  +0.24%  "

# Row 25
$ws.Range("D25").Value = "2.175.71"
$ws.Range("E25").Value = "  +1.07%  "

# Row 26
$origStyle_26 = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.093"
$ws.Range("D26").Style = $origStyle_26
$ws.Range("E26").Value = "  +0.48%  "

# Row 27
$origStyle_27 = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "158.05"
$ws.Range("D27").Style = $origStyle_27

# Row 28
$origStyle_28 = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.40"
$ws.Range("D28").Style = $origStyle_28
$ws.Range("E28").Value = "  +0.35%  "

# Row 29
$origStyle_29 = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.750"
$ws.Range("D29").Style = $origStyle_29
$ws.Range("E29").Value = "  -1.96%  "

# Row 30
$origStyle_30 = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "118.79"
$ws.Range("D30").Style = $origStyle_30
$ws.Range("E30").Value = "  +0.48%  "

# Row 31
$origStyle_31 = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.831"
$ws.Range("D31").Style = $origStyle_31
$ws.Range("E31").Value = "  -0.67%  "

# Row 32
$origStyle_32 = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09324"
$ws.Range("D32").Style = $origStyle_32
$ws.Range("E32").Value = "  +0.34%  "

# Row 33
$origStyle_33 = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.8598"
$ws.Range("D33").Style = $origStyle_33
$ws.Range("E33").Value = "  -2.43%  "

# Row 34
$origStyle_34 = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.157"
$ws.Range("D34").Style = $origStyle_34
$ws.Range("E34").Value = "  -0.27%  "

# Row 35
$origStyle_35 = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.299"
$ws.Range("D35").Style = $origStyle_35
$ws.Range("E35").Value = "  -0.49%  "

# Row 36
$origStyle_36 = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.081"
$ws.Range("D36").Style = $origStyle_36
$ws.Range("E36").Value = "  -1.82%  "

# Row 37
$origStyle_37 = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05781"
$ws.Range("D37").Style = $origStyle_37
$ws.Range("E37").Value = "  +0.55%  "

# Row 38
$origStyle_38 = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.157"
$ws.Range("D38").Style = $origStyle_38
$ws.Range("E38").Value = "  -0.84%  "

# Row 39
$origStyle_39 = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02077"
$ws.Range("D39").Style = $origStyle_39
$ws.Range("E39").Value = "  -0.28%  "

# Row 40
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$origStyle_40 = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5648"
$ws.Range("D40").Style = $origStyle_40
$ws.Range("E40").Value = "  -0.48%  "

# Row 41
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$origStyle_41 = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.616"
$ws.Range("D41").Style = $origStyle_41
$ws.Range("E41").Value = "  -0.46%  "

# Row 42
$origStyle_42 = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1773"
$ws.Range("D42").Style = $origStyle_42
$ws.Range("E42").Value = "  -1.36%  "

# Row 43
$origStyle_43 = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.000002995"
$ws.Range("D43").Style = $origStyle_43
$ws.Range("E43").Value = "  +8.82%  "

# Row 44
$origStyle_44 = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "9.359"
$ws.Range("D44").Style = $origStyle_44
$ws.Range("E44").Value = "  -3.00%  "

# Row 45
$origStyle_45 = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.755"
$ws.Range("D45").Style = $origStyle_45
$ws.Range("E45").Value = "  +8.25%  "

# Row 46
$origStyle_46 = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5268"
$ws.Range("D46").Style = $origStyle_46
$ws.Range("E46").Value = "  -0.87%  "

# Row 47
$origStyle_47 = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "11.38"
$ws.Range("D47").Style = $origStyle_47
$ws.Range("E47").Value = "  -3.50%  "

# Row 48
$origStyle_48 = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06848"
$ws.Range("D48").Style = $origStyle_48
$ws.Range("E48").Value = "  -1.23%  "

# Row 49
$origStyle_49 = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.063"
$ws.Range("D49").Style = $origStyle_49
$ws.Range("E49").Value = "  -5.66%  "

# Row 50
$origStyle_50 = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.806"
$ws.Range("D50").Style = $origStyle_50
$ws.Range("E50").Value = "  -1.19%  "

# Row 51
$origStyle_51 = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "111.47"
$ws.Range("D51").Style = $origStyle_51
$ws.Range("E51").Value = "  -0.60%  "
